$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot the original content (columns A:AY) of every affected
# row before any writes happen, since rows are being permuted amongst
# each other (some in cycles longer than a simple pairwise swap).
$snap = @{}
$snap[8] = $ws.Range("A8:AY8").Value2
$snap[9] = $ws.Range("A9:AY9").Value2
$snap[10] = $ws.Range("A10:AY10").Value2
$snap[11] = $ws.Range("A11:AY11").Value2
$snap[12] = $ws.Range("A12:AY12").Value2
$snap[13] = $ws.Range("A13:AY13").Value2
$snap[14] = $ws.Range("A14:AY14").Value2
$snap[15] = $ws.Range("A15:AY15").Value2
$snap[16] = $ws.Range("A16:AY16").Value2
$snap[17] = $ws.Range("A17:AY17").Value2
$snap[21] = $ws.Range("A21:AY21").Value2
$snap[23] = $ws.Range("A23:AY23").Value2
$snap[25] = $ws.Range("A25:AY25").Value2
$snap[26] = $ws.Range("A26:AY26").Value2
$snap[30] = $ws.Range("A30:AY30").Value2
$snap[31] = $ws.Range("A31:AY31").Value2
$snap[32] = $ws.Range("A32:AY32").Value2
$snap[35] = $ws.Range("A35:AY35").Value2
$snap[36] = $ws.Range("A36:AY36").Value2
$snap[44] = $ws.Range("A44:AY44").Value2
$snap[45] = $ws.Range("A45:AY45").Value2
$snap[46] = $ws.Range("A46:AY46").Value2
$snap[47] = $ws.Range("A47:AY47").Value2
$snap[48] = $ws.Range("A48:AY48").Value2
$snap[49] = $ws.Range("A49:AY49").Value2

# --- Step 2: write each row using the snapshot captured from its designated
# source row (the row whose original content now belongs here).
$ws.Range("A8:AY8").Value2 = $snap[9]
$ws.Range("A9:AY9").Value2 = $snap[8]
$ws.Range("A10:AY10").Value2 = $snap[13]
$ws.Range("A11:AY11").Value2 = $snap[14]
$ws.Range("A12:AY12").Value2 = $snap[15]
$ws.Range("A13:AY13").Value2 = $snap[11]
$ws.Range("A14:AY14").Value2 = $snap[10]
$ws.Range("A15:AY15").Value2 = $snap[12]
$ws.Range("A16:AY16").Value2 = $snap[17]
$ws.Range("A17:AY17").Value2 = $snap[16]
$ws.Range("A21:AY21").Value2 = $snap[23]
$ws.Range("A23:AY23").Value2 = $snap[21]
$ws.Range("A25:AY25").Value2 = $snap[26]
$ws.Range("A26:AY26").Value2 = $snap[25]
$ws.Range("A30:AY30").Value2 = $snap[32]
$ws.Range("A31:AY31").Value2 = $snap[30]
$ws.Range("A32:AY32").Value2 = $snap[31]
$ws.Range("A35:AY35").Value2 = $snap[36]
$ws.Range("A36:AY36").Value2 = $snap[35]
$ws.Range("A44:AY44").Value2 = $snap[49]
$ws.Range("A45:AY45").Value2 = $snap[47]
$ws.Range("A46:AY46").Value2 = $snap[44]
$ws.Range("A47:AY47").Value2 = $snap[46]
$ws.Range("A48:AY48").Value2 = $snap[45]
$ws.Range("A49:AY49").Value2 = $snap[48]

# --- Step 3: Excels COM layer auto-coerces numeric-looking / date-looking
# text into real numbers or date serials when it is read back through
# Value2 and re-assigned. The source cells below are genuinely text, so
# force them back to Text (number format "@") and re-write the literal
# string, then restore "General" formatting to avoid leaving a stray
# explicit text format on the cell.
function Set-TextValue($ws, $addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
}

Set-TextValue $ws "Y8" "2026-01-26"
Set-TextValue $ws "AA8" "2026-01-26"
Set-TextValue $ws "Y9" "2026-01-26"
Set-TextValue $ws "AA9" "2026-01-26"
Set-TextValue $ws "Y10" "2026-01-26"
Set-TextValue $ws "AA10" "2026-01-26"
Set-TextValue $ws "Y11" "2026-01-26"
Set-TextValue $ws "AA11" "2026-01-26"
Set-TextValue $ws "Y12" "2026-01-26"
Set-TextValue $ws "AA12" "2026-01-26"
Set-TextValue $ws "Y13" "2026-01-26"
Set-TextValue $ws "AA13" "2026-01-26"
Set-TextValue $ws "Y14" "2026-01-26"
Set-TextValue $ws "AA14" "2026-01-26"
Set-TextValue $ws "Y15" "2026-01-26"
Set-TextValue $ws "AA15" "2026-01-26"
Set-TextValue $ws "Y16" "2026-01-26"
Set-TextValue $ws "AA16" "2026-01-26"
Set-TextValue $ws "Y17" "2026-01-26"
Set-TextValue $ws "AA17" "2026-01-26"
Set-TextValue $ws "Y21" "2026-01-26"
Set-TextValue $ws "AA21" "2026-01-26"
Set-TextValue $ws "Y23" "2026-01-26"
Set-TextValue $ws "AA23" "2026-01-26"
Set-TextValue $ws "Y25" "2026-01-26"
Set-TextValue $ws "AA25" "2026-01-26"
Set-TextValue $ws "Y26" "2026-01-26"
Set-TextValue $ws "AA26" "2026-01-26"
Set-TextValue $ws "Y30" "2026-01-26"
Set-TextValue $ws "AA30" "2026-01-26"
Set-TextValue $ws "Y31" "2026-01-26"
Set-TextValue $ws "AA31" "2026-01-26"
Set-TextValue $ws "Y32" "2026-01-26"
Set-TextValue $ws "AA32" "2026-01-26"
Set-TextValue $ws "Y35" "2026-01-26"
Set-TextValue $ws "AA35" "2026-01-26"
Set-TextValue $ws "Y36" "2026-01-26"
Set-TextValue $ws "AA36" "2026-01-26"
Set-TextValue $ws "Y44" "2026-01-26"
Set-TextValue $ws "AA44" "2026-01-26"
Set-TextValue $ws "Y45" "2026-01-26"
Set-TextValue $ws "AA45" "2026-01-26"
Set-TextValue $ws "Y46" "2026-01-26"
Set-TextValue $ws "AA46" "2026-01-26"
Set-TextValue $ws "Y47" "2026-01-26"
Set-TextValue $ws "AA47" "2026-01-26"
Set-TextValue $ws "Y48" "2026-01-26"
Set-TextValue $ws "AA48" "2026-01-26"
Set-TextValue $ws "Y49" "2026-01-26"
Set-TextValue $ws "AA49" "2026-01-26"

# Row 26 inherits row 25s original I-column text "1" (a species count
# recorded as text, not a number) — restore it explicitly.
Set-TextValue $ws "I26" "1"

"Row content permutation applied."